$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7175
$ws.Range("C3").Value = 156275
$ws.Range("C4").Value = 147368
$ws.Range("C7").Value = 5.7
$ws.Range("C8").Value = 63.71
